# Apply the "normalization list-01 question-3" edits:
# - B13 ({idPedido, codProduto} ~=> codProduto) becomes ~=> nomeProduto
# - B18 (codProduto -> codProduto) becomes codProduto -> nomeProduto
# - Selection moves from L5 to the merged range B17:E17

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B13").Value = "{idPedido, codProduto} ~=> nomeProduto"
$ws.Range("B18").Value = "codProduto -> nomeProduto"

$ws.Range("B17:E17").Select()
